# Update cryptocurrency price/volume data per Nov 18 2024 refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '91.923.07'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.02%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.121.76'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.88%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '242.72'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +1.66%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '625.91'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -2.59%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.15'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +5.34%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.376'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.00%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '3.120.20'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.44%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.770'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +6.18%  '
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +3.38%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000257'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +2.05%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.60'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -2.92%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '91.724.96'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +1.14%  '
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -1.70%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.703.83'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.78%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.142.01'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.21%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.75'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.65%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0000224'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +3.08%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.75'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +1.94%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.82'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +2.76%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '448.24'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.48%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.63%  '
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -2.21%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '91.89'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.72%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.97'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -3.78%  '
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.00%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.253'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +24.12%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.186'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +15.73%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '9.28'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -4.54%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +34.30%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.01'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +4.90%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +11.03%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '26.68'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -1.54%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +5.54%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.10'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +20.49%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.68'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -5.62%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -1.34%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '493.81'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -4.72%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.99%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.423'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.03%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '22.18'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.00%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.82%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '155.90'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +2.98%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -1.75%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.31%  '
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.60%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '44.82'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -1.92%  '
